$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting existing rows 6-12 down to 7-13
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the latest weekly price entry
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 45203
$ws.Cells.Item(6, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = 100112030
$ws.Cells.Item(6, 7).Value = "Poroto granado"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 800
$ws.Cells.Item(6, 11).Value = 1800
$ws.Cells.Item(6, 12).Value = 2000
$ws.Cells.Item(6, 13).Value = 1900
$ws.Cells.Item(6, 14).Value = "`$/kilo"
$ws.Cells.Item(6, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(6, 16).Value = 1900
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"
